$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.471.20'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.896.69'
$ws.Range('E3').Value = '  +1.01%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.001'
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '238.33'
$ws.Range('E5').Value = '  +1.12%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$ws.Range('E6').Value = '  -0.01%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4918'
$ws.Range('E7').Value = '  +1.00%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2923'
$ws.Range('E8').Value = '  +0.52%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06691'
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('D10').Value = '1.902.25'
$ws.Range('E10').Value = '  +1.22%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '16.91'
$ws.Range('E11').Value = '  +1.89%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.07326'
$ws.Range('E12').Value = '  +1.24%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.173'
$ws.Range('E13').Value = '  +3.60%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '87.50'
$ws.Range('E14').Value = '  -1.58%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.6653'
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '30.450.13'
$ws.Range('E16').Value = '  -0.23%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '13.47'
$ws.Range('E17').Value = '  +4.05%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000007841'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '2.135.99'
$ws.Range('E20').Value = '  +0.74%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.332'
$ws.Range('E21').Value = '  +13.31%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '1.001'
$ws.Range('E22').Value = '  -0.11%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '192.45'
$ws.Range('E23').Value = '  +0.23%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.104'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +2.34%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '162.16'
$ws.Range('E26').Value = '  +2.81%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.22'
$ws.Range('E27').Value = '  -0.59%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.934'
$ws.Range('E28').Value = '  +6.04%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.472'
$ws.Range('E29').Value = '  +4.84%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.329'
$ws.Range('E30').Value = '  +2.15%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.09161'
$ws.Range('E31').Value = '  +1.81%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.048'
$ws.Range('E32').Value = '  +3.28%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.05169'
$ws.Range('E33').Value = '  +0.77%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7384'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  +2.17%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.722'
$ws.Range('E36').Value = '  +1.15%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.01809'
$ws.Range('E37').Value = '  -0.30%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.678'
$ws.Range('E38').Value = '  +0.68%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.9239'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('E40').Value = '  -0.42%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.4385'
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '106.92'
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.910'
$ws.Range('E43').Value = '  +3.64%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.9940'
$ws.Range('E44').Value = '  -0.12%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '68.60'
$ws.Range('E45').Value = '  +20.41%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.1366'
$ws.Range('E46').Value = '  +2.71%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.584'
$ws.Range('E47').Value = '  +3.30%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '8.988'
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('E49').Value = '  +5.62%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05851'
$ws.Range('E50').Value = '  +0.52%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.3915'
$ws.Range('E51').Value = '  -2.35%  '
